# Split the "Bonaire, Sint Eustatius and Saba" row (BES) into three rows,
# one per constituent country, leaving the foreign_tourists values blank
# until they're sourced separately (per commit message: "leaving slots
# for B/E/S in tourism data").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the existing BES row (row 7) so it
# becomes three rows: 7 (Bonaire), 8 (Sint Eustatius), 9 (Saba).
$ws.Rows("8:9").Insert()

# Row 8: Sint Eustatius (name/name_govt stays the umbrella name)
$ws.Cells.Item(8, 1).Value = "Bonaire, Sint Eustatius and Saba"
$ws.Cells.Item(8, 2).Style = "Normal"
$ws.Cells.Item(8, 2).Value = "Sint Eustatius"
$ws.Cells.Item(8, 4).Value = 535
$ws.Cells.Item(8, 5).Clear()

# Row 9: Saba
$ws.Cells.Item(9, 1).Value = "Bonaire, Sint Eustatius and Saba"
$ws.Cells.Item(9, 2).Style = "Normal"
$ws.Cells.Item(9, 2).Value = "Saba"
$ws.Cells.Item(9, 4).Value = 535
$ws.Cells.Item(9, 5).Clear()

# Row 7: retarget the original row to just Bonaire and clear its old
# foreign_tourists figure (it applied to the whole BES group).
$ws.Cells.Item(7, 2).Value = "Bonaire"
$ws.Cells.Item(7, 3).Value = "BESB"
$ws.Cells.Item(7, 5).Clear()

# Fill in the alpha_3-style codes for the two new rows.
$ws.Cells.Item(8, 3).Value = "BESE"
$ws.Cells.Item(9, 3).Value = "BESS"

# Update the view: selection moves to E7, no more frozen/scrolled
# top-left cell.
$ws.Range("E7").Select() | Out-Null
